# Adjust topic schedule to move API access next to web scraping.
#
# The "Web scraping" and "API access" topics (previously on separate,
# non-adjacent weeks) are combined into a single two-week "Getting data
# from the web" unit. Every topic scheduled after them shifts up one
# slot to fill the gap left by removing the standalone "API access"
# week, and "Introduction to Python"/"Functional programming in Python"
# each gain an extra week in the process.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D11").Value = "Introduction to Python"
$ws.Range("D12").Value = "Functional programming in Python"
$ws.Range("D13").Value = "Getting data from the web: API access"
$ws.Range("D14").Value = "Getting data from the web: scraping"
$ws.Range("D15").Value = "Network analysis"
$ws.Range("D16").Value = "Text analysis"
$ws.Range("D17").Value = "Machine learning"
$ws.Range("D18").Value = "Data warehousing and relational databases"

# Move the active selection to D15, matching where the author's cursor
# ended up after making the edit.
[void]$ws.Range("D15").Select()
